$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 and 3 (A and B values change; C/D/E unchanged) ---
$ws.Range("A2").Value = 2148226
$ws.Range("B2").Value = 78181

$ws.Range("A3").Value = 2159417
$ws.Range("B3").Value = 58395

# --- Extend formatting for new rows 4-12 by copying row 3's formats ---
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new data rows 4-12 ---
$data = @(
    @(4, 2159872, 112094),
    @(5, 2164252, 41551),
    @(6, 2168432, 25581),
    @(7, 2177621, 87923),
    @(8, 2181725, 30916),
    @(9, 2181750, 28931),
    @(10, 2185057, 73218),
    @(11, 2188687, 49081),
    @(12, 2188702, 125119)
)

foreach ($row in $data) {
    $r = $row[0]
    $orderVal = $row[1]
    $sumVal = $row[2]
    $ws.Cells.Item($r, 1).Value = $orderVal
    $ws.Cells.Item($r, 2).Value = $sumVal
    $ws.Cells.Item($r, 3).Value = 7
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 3.5
}

# --- Update the selected cell to match the saved view state ---
$ws.Range("B17").Select()

Write-Output "done"
